# Apply the edits described by the diff to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric "Qty executed upto date" (column C) updates ---
$ws.Range("C8").Value = 99
$ws.Range("C9").Value = 14
$ws.Range("C10").Value = 17
$ws.Range("C11").Value = 26
$ws.Range("C12").Value = 99
$ws.Range("C13").Value = 96
$ws.Range("C15").Value = 50
$ws.Range("C16").Value = 53
$ws.Range("C17").Value = 12

# --- Text-formatted "Amount" columns (G/H) updates ---
# These cells store their values as plain text (e.g. "3584.00"), so force
# a text number format before assigning, otherwise Excel would coerce the
# numeric-looking string back into a number and drop the text semantics.
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "3584.00"

$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "8024.00"

$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "17212.00"

$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "13056.00"

$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "42566.00"

$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "42566.00"

$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "42566.00"

$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "42566.00"
